$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row formatting (bold/centered/top/bordered) to the newly added columns L:O
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:O1").PasteSpecial(-4122) | Out-Null

# Header row text
$ws.Range("A1").Value = "Best Estimator"
$ws.Range("B1").Value = "Best Score"
$ws.Range("C1").Value = "Best Params"
$ws.Range("D1").Value = "CV Train F1"
$ws.Range("E1").Value = "CV Test F1"
$ws.Range("F1").Value = "Validation F1"
$ws.Range("G1").Value = "CV Train Precision"
$ws.Range("H1").Value = "CV Test Precision"
$ws.Range("I1").Value = "Validation Precision"
$ws.Range("J1").Value = "CV Train Recall"
$ws.Range("K1").Value = "CV Test Recall"
$ws.Range("L1").Value = "Validation Recall"
$ws.Range("M1").Value = "Y Val (Validation)"
$ws.Range("N1").Value = "Y Pred (Validation)"
$ws.Range("O1").Value = "Seed"

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),`n                ('model', SVC(C=5, class_weight='balanced', random_state=42))])"
$ws.Range("B2").Value = 0.7733333333333333
$ws.Range("C2").Value = "{'selector': None, 'scaler': RobustScaler(), 'model__kernel': 'rbf', 'model__class_weight': 'balanced', 'model__C': 5}"
$ws.Range("D2").Value = 0.6795367153943095
$ws.Range("E2").Value = 0.5364411167536167
$ws.Range("F2").Value = 0.6382978723404256
$ws.Range("G2").Value = 0.6320134184083219
$ws.Range("H2").Value = 0.4932014715608465
$ws.Range("I2").Value = 0.6
$ws.Range("J2").Value = 0.7798313492063492
$ws.Range("K2").Value = 0.6450833333333333
$ws.Range("L2").Value = 0.6818181818181818
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1 0 1 1 1 1 0 0 0 0 1 0 1]"
$ws.Range("N2").Value = "[1 0 1 1 1 0 1 0 1 1 0 1 1 1 1 0 1 0 0 1 1 1 0 0 1 0 1 1 1 1 1 1 0 1 1 1]"
$ws.Range("O2").Value = 42

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None), ('selector', None),`n                ('model', SVC(C=5, kernel='sigmoid', random_state=42))])"
$ws.Range("B3").Value = 0.7207142857142858
$ws.Range("C3").Value = "{'selector': None, 'scaler': None, 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 5}"
$ws.Range("D3").Value = 0.6260414013189817
$ws.Range("E3").Value = 0.4867197293447293
$ws.Range("F3").Value = 0.8076923076923077
$ws.Range("G3").Value = 0.6057127092645993
$ws.Range("H3").Value = 0.4497985284391534
$ws.Range("I3").Value = 0.75
$ws.Range("J3").Value = 0.6953437499999999
$ws.Range("K3").Value = 0.5835416666666667
$ws.Range("L3").Value = 0.875
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1]"
$ws.Range("N3").Value = "[1 0 1 1 1 1 1 1 1 1 1 0 0 1 1 1 1 1 1 1 0 1 1 0 0 1 1 1 1 0 0 1 1 1 1 1]"
$ws.Range("O3").Value = 69

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None), ('selector', None),`n                ('model',`n                 SVC(C=0.0001, class_weight='balanced', kernel='linear',`n                     random_state=42))])"
$ws.Range("B4").Value = 0.7521428571428571
$ws.Range("C4").Value = "{'selector': None, 'scaler': None, 'model__kernel': 'linear', 'model__class_weight': 'balanced', 'model__C': 0.0001}"
$ws.Range("D4").Value = 0.6699718221911494
$ws.Range("E4").Value = 0.517059857966108
$ws.Range("F4").Value = 0.5777777777777778
$ws.Range("G4").Value = 0.6547855553499173
$ws.Range("H4").Value = 0.4891123511904762
$ws.Range("I4").Value = 0.6842105263157895
$ws.Range("J4").Value = 0.7490021929824562
$ws.Range("K4").Value = 0.6024166666666667
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 0 1 0 1 1 1 0 1]"
$ws.Range("N4").Value = "[1 1 1 0 0 0 0 1 1 0 1 0 0 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 1 0 1 1 1 1 0]"
$ws.Range("O4").Value = 23

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),`n                ('model', SVC(C=1, class_weight='balanced', random_state=42))])"
$ws.Range("B5").Value = 0.7285714285714285
$ws.Range("C5").Value = "{'selector': None, 'scaler': RobustScaler(), 'model__kernel': 'rbf', 'model__class_weight': 'balanced', 'model__C': 1}"
$ws.Range("D5").Value = 0.6776677790535152
$ws.Range("E5").Value = 0.5443849194786695
$ws.Range("F5").Value = 0.6666666666666667
$ws.Range("G5").Value = 0.6289911268915526
$ws.Range("H5").Value = 0.4799864417989418
$ws.Range("I5").Value = 0.7647058823529411
$ws.Range("J5").Value = 0.7810962301587302
$ws.Range("K5").Value = 0.675
$ws.Range("L5").Value = 0.5909090909090909
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 1 1 1 1 1 0]"
$ws.Range("N5").Value = "[0 0 1 0 1 0 0 0 0 0 1 1 0 0 1 0 0 0 1 1 1 1 1 1 1 0 0 1 0 1 0 1 1 1 0 0]"
$ws.Range("O5").Value = 99

# Row 6
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',`n                                                     random_state=42))),`n                ('model', SVC(C=5, kernel='sigmoid', random_state=42))])"
$ws.Range("B6").Value = 0.7239285714285714
$ws.Range("C6").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 5}"
$ws.Range("D6").Value = 0.7936093157032132
$ws.Range("E6").Value = 0.6437831277518777
$ws.Range("F6").Value = 0.608695652173913
$ws.Range("G6").Value = 0.725395209189486
$ws.Range("H6").Value = 0.5521510416666666
$ws.Range("I6").Value = 0.5384615384615384
$ws.Range("J6").Value = 0.9106770833333333
$ws.Range("K6").Value = 0.8117083333333334
$ws.Range("L6").Value = 0.7
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1 1 0 1 0 1 1 1 1 1 1 1 0]"
$ws.Range("N6").Value = "[1 1 1 0 1 1 1 1 0 1 1 1 1 1 1 0 1 1 0 1 0 1 1 1 1 0 1 1 0 1 1 0 1 0 0 1]"
$ws.Range("O6").Value = 89
